$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column-width metadata housekeeping -----------------------------------
# The lone "customWidth" marker that originally sat on column Q (17) shifts to
# column R (18) once the Long Run table grows into column Q. Inserting (and
# immediately it pushes the pre-existing per-column
# formatting one slot to the right without disturbing any cell data, since at
# this point nothing has been written past column K yet.
$ws.Columns("Q").Insert()

# --- New header row (row 10/11) --------------------------------------------
$ws.Range("Q11").Value = "Long Run"
$ws.Range("K10").Value = "Mileage"
$ws.Range("K11").Value = "Planned"
$ws.Range("J11").Value = "Equation"
$ws.Range("L11").Value = "Percent"

# --- Move the Long Run formulas from column K to column Q ------------------
$ws.Range("Q12").Formula = "=Q13-2"
$ws.Range("Q13").Formula = "=Q14-2"
$ws.Range("Q14").Formula = "=Q15-2"
$ws.Range("Q15").Formula = "=Q16-2"
$ws.Range("Q16").Formula = "=Q18-1"
$ws.Range("Q17").Formula = "=Q19"
$ws.Range("Q18").Formula = "=Q20-1"
$ws.Range("Q19").Formula = "=Q21-1"
$ws.Range("Q20").Formula = "=maxLR"
$ws.Range("Q21").Formula = "=maxLR * 0.85"
$ws.Range("Q22").Formula = "=maxLR"
$ws.Range("Q23").Formula = "=maxLR * 0.85"
$ws.Range("Q24").Formula = "=maxLR"
$ws.Range("Q25").Formula = "=maxLR * 0.8"
$ws.Range("Q26").Formula = "=maxLR*0.6"
$ws.Range("Q27").Formula = "=goalDistance"

# --- New "Planned" mileage numbers in column K ------------------------------
$plannedMileage = @{
    12 = 25
    13 = 28
    14 = 33
    15 = 37
    16 = 40
    17 = 44
    18 = 48
    19 = 40
    20 = 50
    21 = 42
    22 = 50
    23 = 42
    24 = 50
    25 = 40
    26 = 30
    27 = 20
}
foreach ($row in 12..27) {
    $ws.Range("K$row").Value = $plannedMileage[$row]
}

# --- New "Equation" column: blank except for peak weeks that equal maxMileage
$ws.Range("J12").Value = ""
$ws.Range("J13").Value = ""
$ws.Range("J14").Value = ""
$ws.Range("J15").Value = ""
$ws.Range("J16").Value = ""
$ws.Range("J17").Value = ""
$ws.Range("J18").Value = ""
$ws.Range("J19").Value = ""
$ws.Range("J20").Formula = "=maxMileage"
$ws.Range("J21").Value = ""
$ws.Range("J22").Formula = "=maxMileage"
$ws.Range("J23").Value = ""
$ws.Range("J24").Formula = "=maxMileage"
$ws.Range("J25").Value = ""
$ws.Range("J26").Value = ""
$ws.Range("J27").Value = ""

# --- New "Percent" column -----------------------------------------------
foreach ($row in 12..27) {
    $ws.Range("L$row").Formula = "=K$row/maxMileage"
}

# --- Alignment / styling ----------------------------------------------------
# Horizontal-center goes onto the data grid (creates style index 2)...
$ws.Range("J12:K27").HorizontalAlignment = -4108
# ...vertical-center goes onto the three new headers (creates style index 3)
$ws.Range("J11:L11").VerticalAlignment = -4108

# --- Selection / view state --------------------------------------------------
$ws.Range("U18:W24").Select()

$wb.Windows.Item(1).WindowState = -4137
